# New docking path: add a recomputed set of docking points (Startpunkt /
# Zwischenwert / Einfahrt / Endpunkt) together with a new "shift 10 cm to
# the left" offset row that those points are derived from.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: Startpunkt shifted by the new offset.
$ws.Range("B16").Value = "Startpunkt"
$ws.Range("C16").Formula = "=C5+C`$14"
$ws.Range("D16").Formula = "=D5+D`$14"

# Row 17: Zwischenwert shifted by the new offset.
$ws.Range("B17").Value = "Zwischenwert"
$ws.Range("C17").Formula = "=C6+C`$14"
$ws.Range("D17").Formula = "=D6+D`$14"

# Row 18: Einfahrt shifted by the new offset (A18 already holds dx0 text).
$ws.Range("B18").Value = "Einfahrt"
$ws.Range("C18").Formula = "=C7+C`$14"
$ws.Range("D18").Formula = "=D7+D`$14"

# Row 19: Endpunkt shifted by the new offset (A19 already holds dy0 text).
$ws.Range("B19").Value = "Endpunkt"
$ws.Range("C19").Formula = "=C8+C`$14"
$ws.Range("D19").Formula = "=D8+D`$14"

# Row 14: new "Verschiebung nach links 10 cm" (shift left 10 cm) offset
# that rows 16-19 above reference via C$14 / D$14.
$ws.Range("B14").Value = "Verschiebung nach links 10 cm"
$ws.Range("C14").Formula = "=-D12*E14"
$ws.Range("D14").Formula = "=C12*E14"
$ws.Range("E14").Value = 0.1

# Selection moves to F26 (next free area below the new table).
$ws.Range("F26").Select() | Out-Null
